$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the course content template: "Semester 14" -> "Semester 1"
$ws.Range("A2").Value = "Semester 1"

# Move the active selection to A3, as Excel does after committing an edit
$ws.Range("A3").Select()
